$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header labels (shared strings used in C1 and D1)
$ws.Range("C1").Value = "CV Mean Accuracy"
$ws.Range("D1").Value = "CV Std Dev"

# Update numeric values for KNN / Depression row (row 6)
$ws.Range("C6").Value = 0.8375
$ws.Range("D6").Value = 0.03061862178478972

# Update numeric values for KNN / Anxiety row (row 7)
$ws.Range("C7").Value = 0.8625
$ws.Range("D7").Value = 0.04677071733467426

# Update numeric values for KNN / Panic Attack row (row 10)
$ws.Range("C10").Value = 0.6125
$ws.Range("D10").Value = 0.2833946012188658
